$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.444.44'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '1.803.05'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '227.58'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("E6").Value = '  +3.28%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '36.31'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +7.18%  '
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("E10").Value = '  -0.35%  '
$ws.Range("E11").Value = '  +1.61%  '
$ws.Range("D12").Value = '2.062.41'
$ws.Range("E12").Value = '  -0.68%  '
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '11.57'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.02%  '
$ws.Range("D14").Value = '1.800.59'
$ws.Range("E14").Value = '  -1.04%  '
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '4.49'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +4.55%  '
$ws.Range("D17").Value = '34.398.47'
$ws.Range("E17").Value = '  -0.20%  '
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '69.47'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '245.10'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("E20").Value = '  -1.20%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '11.59'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.46%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '2.16'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.15%  '
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '172.50'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.69%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '7.98'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +8.54%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '16.86'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.68%  '
$ws.Range("E28").Value = '  +1.63%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  +0.62%  '
$ws.Range("E31").Value = '  -0.12%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '3.83'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("E33").Value = '  -0.28%  '
$ws.Range("E34").Value = '  -1.81%  '
$ws.Range("D35").Value = '1.396.72'
$ws.Range("E35").Value = '  -0.91%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.674'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.48%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '2.43'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.62%  '
$ws.Range("E38").Value = '  -0.62%  '
$ws.Range("E39").Value = '  -0.34%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '82.45'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.03%  '
$ws.Range("E41").Value = '  +0.34%  '
$ws.Range("E42").Value = '  -0.87%  '
$ws.Range("E43").Value = '  +0.41%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '1.19'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +6.85%  '
$ws.Range("E45").Value = '  -5.51%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '6.04'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.35%  '
$ws.Range("E47").Value = '  -3.96%  '
$ws.Range("D48").Value = '1.963.67'
$ws.Range("E48").Value = '  -0.63%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '104.47'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.26%  '
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("D51").Value = '0.0₆0127'
$ws.Range("E51").Value = '  -0.54%  '
